# Generate Report for handback
# Update the "Correspond Handoff Datetime" (D2) and "Correspond Handback DateTime" (G2)
# timestamps for the first data row on the "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-11 03:02:54"
$wsZhCn.Range("G2").Value = "2016-01-11 03:04:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-11 03:03:10"
$wsDeDe.Range("G2").Value = "2016-01-11 03:04:28"
